# Applies the RTM bug/test-case remapping edit described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header typo fix -------------------------------------------------
$ws.Range("I5").Value = "Current Status"

# --- Row 6 : Test Case range unchanged, Bug list collapses to Bug 001-006
$ws.Range("J6").Value = "Bug - 001 / Bug - 002 / Bug - 003 / Bug - 004 / Bug - 005 / Bug - 006"

# --- Row 7 -------------------------------------------------------------
$ws.Range("H7").Value = "TC - 032 to TC - 039"
$ws.Range("J7").Value = "Bug - 007"

# --- Row 8 -------------------------------------------------------------
$ws.Range("H8").Value = "TC - 040 to TC - 044"
$ws.Range("J8").Value = "Bug - 008"

# --- Row 9 -------------------------------------------------------------
$ws.Range("H9").Value = "TC - 045 to TC - 074"
$ws.Range("J9").Value = "Bug - 009 / Bug - 010 / Bug - 011 / Bug - 012 / Bug - 013 / Bug - 014"

# --- Row 10 ------------------------------------------------------------
$ws.Range("H10").Value = "TC - 075 to TC - 102"
$ws.Range("J10").Value = "Bug - 015"

# --- Row 11 ------------------------------------------------------------
$ws.Range("H11").Value = "TC - 103 to TC - 117"
$ws.Range("J11").Value = "Bug - 016 / Bug - 017 / Bug - 018 / Bug - 019 / Bug - 020"

# --- Row 12 (status flips Failed -> Passed, bug list cleared) ----------
$ws.Range("H12").Value = "TC - 118 to TC - 129"
$ws.Range("I15").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I12").Value = "Passed"
$ws.Range("J12").Value = "-"

# --- Row 13 --------------------------------------------------------------
$ws.Range("H13").Value = "TC - 130 to TC - 135"
$ws.Range("J13").Value = "Bug - 021 / Bug - 022 / Bug - 023"

# --- Row 14 (status flips Failed -> Passed, bug list cleared) ----------
$ws.Range("H14").Value = "TC - 136 to TC - 145"
$ws.Range("I15").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I14").Value = "Passed"
$ws.Range("J14").Value = "-"

# --- Row 15 --------------------------------------------------------------
$ws.Range("H15").Value = "TC - 146 to TC - 154"
$ws.Range("J15").Value = "-"

# --- Row 16 --------------------------------------------------------------
$ws.Range("H16").Value = "TC - 155 to TC - 163"
$ws.Range("J16").Value = "Bug - 024"

# --- Row 17 --------------------------------------------------------------
$ws.Range("H17").Value = "TC - 164 to TC - 172"
$ws.Range("J17").Value = "Bug - 025 / Bug - 026 / Bug - 027 / Bug - 028"

# --- Row heights adjust now that text no longer wraps as much ----------
$ws.Rows.Item(6).RowHeight = 43.5
$ws.Rows.Item(9).RowHeight = 43.5

# --- View: scroll down a bit and move the active selection --------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P11").Select() | Out-Null
